# Apply "reinforcement learning is added" edit:
#  - Replace training rows 2-26 with the new labelled data (text / label / value)
#  - Append new rows 27-37 with additional reinforcement-learning examples
#  - Row 19's text flips from "дороу"(idx7) to "привет"(idx3) (moved below)
#  - "привет спасибо" sample is dropped, "мис привет" samples become
#    leading-space variants " привет" / " сказать привет" / " дар"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("дар",              "Не приветствие", 0),
  @("дар",              "Не приветствие", 0),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("холодно",          "Не приветствие", 0),
  @("дороу",            "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @("привет",           "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" сказать привет",  "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" дар",             "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" привет",          "Приветствие",    1),
  @(" дар",             "Не приветствие", 0),
  @("привет",           "Приветствие",    1)
)

$row = 2
foreach ($item in $data) {
  $ws.Cells.Item($row, 1).Value = $item[0]
  $ws.Cells.Item($row, 2).Value = $item[1]
  $ws.Cells.Item($row, 3).Value = $item[2]
  $row = $row + 1
}
